$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.035.97'
$ws.Range('E2').Value = '  +0.39%  '
$ws.Range('D3').Value = '1.644.42'
$ws.Range('E3').Value = '  +0.39%  '
$ws.Range('E4').Value = '  +0.55%  '
$ws.Range('D5').Value = '216.14'
$ws.Range('E5').Value = '  +0.72%  '
$ws.Range('E6').Value = '  +0.30%  '
$ws.Range('E7').Value = '  +0.51%  '
$ws.Range('B8').Value = 'Dogecoin'
$ws.Range('C8').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D8').Value = '0.0639'
$ws.Range('E8').Value = '  +0.50%  '
$ws.Range('B9').Value = 'Cardano'
$ws.Range('C9').Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range('D9').Value = '0.255'
$ws.Range('E9').Value = '  +0.36%  '
$ws.Range('D10').Value = '19.57'
$ws.Range('E10').Value = '  -0.14%  '
$ws.Range('E11').Value = '  +0.44%  '
$ws.Range('E12').Value = '  +0.70%  '
$ws.Range('D13').Value = '1.656.66'
$ws.Range('E13').Value = '  +0.84%  '
$ws.Range('D14').Value = '0.545'
$ws.Range('E14').Value = '  +0.28%  '
$ws.Range('B15').Value = 'ShibaInu'
$ws.Range('C15').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D15').Value = '0.0₃0763'
$ws.Range('E15').Value = '  +0.91%  '
$ws.Range('B16').Value = 'Litecoin'
$ws.Range('C16').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D16').Value = '63.46'
$ws.Range('E16').Value = '  +1.43%  '
$ws.Range('D17').Value = '26.060.62'
$ws.Range('E17').Value = '  +0.43%  '
$ws.Range('E18').Value = '  +0.48%  '
$ws.Range('D19').Value = '194.43'
$ws.Range('E19').Value = '  +0.31%  '
$ws.Range('E20').Value = '  -0.68%  '
$ws.Range('D21').Value = '9.94'
$ws.Range('E21').Value = '  +0.27%  '
$ws.Range('E22').Value = '  -1.11%  '
$ws.Range('E23').Value = '  +4.27%  '
$ws.Range('E24').Value = '  -0.94%  '
$ws.Range('E25').Value = '  +0.50%  '
$ws.Range('D26').Value = '143.26'
$ws.Range('E26').Value = '  -0.36%  '
$ws.Range('E27').Value = '  +0.52%  '
$ws.Range('D28').Value = '15.52'
$ws.Range('E28').Value = '  +0.54%  '
$ws.Range('E29').Value = '  +0.62%  '
$ws.Range('E30').Value = '  -1.12%  '
$ws.Range('E31').Value = '  -0.29%  '
$ws.Range('E32').Value = '  +1.40%  '
$ws.Range('E33').Value = '  -0.65%  '
$ws.Range('E34').Value = '  +1.66%  '
$ws.Range('E35').Value = '  +0.25%  '
$ws.Range('D36').Value = '1.131.15'
$ws.Range('E36').Value = '  -0.64%  '
$ws.Range('E37').Value = '  -0.82%  '
$ws.Range('E38').Value = '  +0.38%  '
$ws.Range('E39').Value = '  +0.09%  '
$ws.Range('E40').Value = '  +0.83%  '
$ws.Range('D41').Value = '99.18'
$ws.Range('D42').Value = '0.798'
$ws.Range('E42').Value = '  -0.03%  '
$ws.Range('E43').Value = '  +2.35%  '
$ws.Range('D44').Value = '56.46'
$ws.Range('E44').Value = '  -0.06%  '
$ws.Range('E46').Value = '  +2.77%  '
$ws.Range('D47').Value = '7.79'
$ws.Range('E47').Value = '  +2.09%  '
$ws.Range('E48').Value = '  -0.05%  '
$ws.Range('E49').Value = '  +0.37%  '
$ws.Range('E50').Value = '  -0.86%  '
$ws.Range('D51').Value = '1.19'
$ws.Range('E51').Value = '  +3.13%  '
